$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 80 (pushes existing rows 80-188 down to 81-189)
$ws.Rows("80:80").Insert()

# Populate the newly inserted row 80 with the new data record.
# Columns A,B,C,E,F,G,H,I,N,Q,R mirror the surrounding rows (constant for this sheet/series).
$ws.Range("A80").Value = 5
$ws.Range("B80").Value = "Macroferia Regional de Talca"
$ws.Range("C80").Value = "Maule"
$ws.Range("D80").Value = 44467
$ws.Range("E80").Value = 7
$ws.Range("F80").Value = 100114013
$ws.Range("G80").Value = "Zanahoria"
$ws.Range("H80").Value = "Sin especificar"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 400
$ws.Range("K80").Value = 6000
$ws.Range("L80").Value = 6000
$ws.Range("M80").Value = 6000
$ws.Range("N80").Value = "`$/saco 20 kilos"
$ws.Range("O80").Value = "Región de Ñuble"
$ws.Range("P80").Value = 300
$ws.Range("Q80").Value = 20
$ws.Range("R80").Value = "Hortaliza"
